$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Prise de RDV ELW" link cells (D8 and D24) to use the
# microsoft-edge: prefixed address, and add the "ignore" marker in column E.
$ws.Range("D8").Value = "microsoft-edge:https://b2b.kiabi.fr/logrdv"
$ws.Range("E8").Value = "ignore"

# D24 previously held a real hyperlink; remove it (it becomes plain text)
# before overwriting the cell value.
$ws.Hyperlinks.Item(2).Delete()
$ws.Range("D24").Value = "microsoft-edge:https://b2b.kiabi.fr/logrdv"
$ws.Range("E24").Value = "ignore"

# Match the author's final selection in the sheet.
$ws.Range("D24:E24").Select()
